# docs: update docs - es/es-oral.xlsx: swap the French question words for
# their Spanish equivalents (the sheet pairs an English question word with
# its translation) and move the sheet's active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "English / French" (x2) -> "English / Spanish" (x2)
$ws.Range("A1").Value = "Spanish"
$ws.Range("B1").Value = "English"
$ws.Range("C1").Value = "Spanish"
$ws.Range("D1").Value = "English"

# French word -> Spanish word (the English translation column is unchanged)
# comment -> cómo (how\what)
$ws.Range("A2").Value = "cómo"
# quand -> cuando (when)
$ws.Range("C2").Value = "cuando"

# quel/quelle -> cual (what\how)
$ws.Range("A3").Value = "cual"
# qui -> quién (who)
$ws.Range("C3").Value = "quién"

# où -> dónde (where)
$ws.Range("A4").Value = "dónde"
# pourquoi -> por qué (why)
$ws.Range("C4").Value = "por qué"

# The saved selection moved from B7 to A5.
$ws.Range("A5").Select()

# Cosmetic: the saved window got slightly wider (29040 -> 29080 twips).
$excel.ActiveWindow.Width = 29080
